$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.240.38'
$ws.Range("E2").Value = '  -0.91%  '

$ws.Range("D3").Value = '2.366.67'
$ws.Range("E3").Value = '  -0.86%  '

$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.12%  '

$ws.Range("D5").Value = '''505.33'
$ws.Range("E5").Value = '  +0.37%  '

$ws.Range("D6").Value = '''129.89'
$ws.Range("E6").Value = '  -1.25%  '

$ws.Range("D7").Value = '''1.01'
$ws.Range("E7").Value = '  +0.82%  '

$ws.Range("D8").Value = '''0.543'
$ws.Range("E8").Value = '  -1.53%  '

$ws.Range("D9").Value = '2.373.94'
$ws.Range("E9").Value = '  -0.77%  '

$ws.Range("D10").Value = '''0.0986'
$ws.Range("E10").Value = '  +1.82%  '

$ws.Range("E11").Value = '  -0.14%  '

$ws.Range("D12").Value = '''4.88'
$ws.Range("E12").Value = '  +7.30%  '

$ws.Range("D13").Value = '''0.325'
$ws.Range("E13").Value = '  +1.39%  '

$ws.Range("D14").Value = '2.788.01'
$ws.Range("E14").Value = '  -0.74%  '

$ws.Range("D15").Value = '56.117.50'
$ws.Range("E15").Value = '  -0.90%  '

$ws.Range("D16").Value = '''21.62'
$ws.Range("E16").Value = '  -0.37%  '

$ws.Range("E17").Value = '  -0.10%  '

$ws.Range("D18").Value = '2.322.68'
$ws.Range("E18").Value = '  -4.67%  '

$ws.Range("D19").Value = '''10.00'
$ws.Range("E19").Value = '  -1.65%  '

$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").Value = '''309.56'
$ws.Range("E20").Value = '  +0.39%  '

$ws.Range("B21").Value = 'Polkadot'
$ws.Range("C21").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D21").Value = '''4.03'
$ws.Range("E21").Value = '  +0.04%  '

$ws.Range("D22").Value = '''6.26'
$ws.Range("E22").Value = '  -0.14%  '

$ws.Range("D23").Value = '''0.999'
$ws.Range("E23").Value = '  -0.04%  '

$ws.Range("D24").Value = '''65.65'
$ws.Range("E24").Value = '  +1.19%  '

$ws.Range("E25").Value = '  -0.48%  '

$ws.Range("D26").Value = '''0.371'
$ws.Range("E26").Value = '  -0.68%  '

$ws.Range("E27").Value = '  -1.63%  '

$ws.Range("D28").Value = '''7.18'
$ws.Range("E28").Value = '  -3.11%  '

$ws.Range("D29").Value = '''173.01'
$ws.Range("E29").Value = '  -1.00%  '

$ws.Range("E30").Value = '  -0.93%  '

$ws.Range("E31").Value = '  -1.16%  '

$ws.Range("D32").Value = '''5.84'
$ws.Range("E32").Value = '  -1.12%  '

$ws.Range("E33").Value = '  +0.01%  '

$ws.Range("E34").Value = '  +0.01%  '

$ws.Range("E35").Value = '  -3.76%  '

$ws.Range("D36").Value = '''17.56'
$ws.Range("E36").Value = '  -1.76%  '

$ws.Range("E37").Value = '  -0.36%  '

$ws.Range("D38").Value = '''3.68'
$ws.Range("E38").Value = '  -3.25%  '

$ws.Range("D39").Value = '''0.824'
$ws.Range("E39").Value = '  +1.44%  '

$ws.Range("D40").Value = '''36.28'
$ws.Range("E40").Value = '  -1.50%  '

$ws.Range("E41").Value = '  -3.18%  '

$ws.Range("D42").Value = '''3.37'
$ws.Range("E42").Value = '  +0.68%  '

$ws.Range("D43").Value = '''125.10'
$ws.Range("E43").Value = '  -5.33%  '

$ws.Range("E44").Value = '  -2.72%  '

$ws.Range("E45").Value = '  -0.07%  '

$ws.Range("D46").Value = '''0.0899'
$ws.Range("E46").Value = '  -0.87%  '

$ws.Range("D47").Value = '''237.03'
$ws.Range("E47").Value = '  -4.81%  '

$ws.Range("E48").Value = '  -0.87%  '

$ws.Range("D49").Value = '''0.0207'
$ws.Range("E49").Value = '  -1.38%  '

$ws.Range("D50").Value = '''16.93'
$ws.Range("E50").Value = '  -0.72%  '

$ws.Range("E51").Value = '  +0.34%  '
